# Rename the header columns to the new (space-free) labels.
# Process Date      -> ProcessDate
# Client Account    -> ClientAccount
# External Reference-> ExternalReference
# Company Name      -> CompanyName
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ProcessDate"
$ws.Range("B1").Value = "ClientAccount"
$ws.Range("C1").Value = "ExternalReference"
$ws.Range("D1").Value = "CompanyName"

# Move the active selection from G6 to G10 (bottom pane, below the frozen header row).
$ws.Range("G10").Select() | Out-Null
